$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "x"
$ws.Range("D4").Value = "y"
$ws.Range("A4").Value = "e"
$ws.Range("B3").Value = "f"
$ws.Range("D5").Value = ""

$ws.Range("B3").Select()
